$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Total" column header
$ws.Range("S1").Value = "Total"

# New per-row totals (column S) for existing rows 2-6
$ws.Range("S2").Value = 2032
$ws.Range("S3").Value = 215
$ws.Range("S4").Value = 1078
$ws.Range("S5").Value = 258
$ws.Range("S6").Value = 1278

# New row 7: "Outros"
$ws.Range("A7").Value = "Outros"
$ws.Range("B7").Value = 157
$ws.Range("C7").Value = 10
$ws.Range("D7").Value = 9
$ws.Range("E7").Value = 64
$ws.Range("F7").Value = 75
$ws.Range("G7").Value = 88
$ws.Range("H7").Value = 91
$ws.Range("I7").Value = 93
$ws.Range("J7").Value = 87
$ws.Range("K7").Value = 116
$ws.Range("L7").Value = 114
$ws.Range("M7").Value = 96
$ws.Range("N7").Value = 93
$ws.Range("O7").Value = 97
$ws.Range("P7").Value = 116
$ws.Range("Q7").Value = 137
$ws.Range("R7").Value = 422
$ws.Range("S7").Value = 1865

# New row 8: "Total"
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Value = 178
$ws.Range("C8").Value = 15
$ws.Range("D8").Value = 14
$ws.Range("E8").Value = 76
$ws.Range("F8").Value = 88
$ws.Range("G8").Value = 116
$ws.Range("H8").Value = 127
$ws.Range("I8").Value = 151
$ws.Range("J8").Value = 181
$ws.Range("K8").Value = 283
$ws.Range("L8").Value = 368
$ws.Range("M8").Value = 442
$ws.Range("N8").Value = 481
$ws.Range("O8").Value = 526
$ws.Range("P8").Value = 683
$ws.Range("Q8").Value = 793
$ws.Range("R8").Value = 2204
$ws.Range("S8").Value = 6726
